$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# --- Crime statistics table updates (rows 15-31) ---

# Row 15
$ws.Range("F15").Value = 6
$ws.Range("I15").Value = 28
$ws.Range("K15").Value = 100
$ws.Range("L15").Value = 180
$ws.Range("M15").Value = 154.545454545455
$ws.Range("N15").Value = -3.448275862068

# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -60
$ws.Range("I16").Value = 68
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -28.421052631578
$ws.Range("L16").Value = -5.555555555555
$ws.Range("M16").Value = -16.049382716049
$ws.Range("N16").Value = -84.821428571428

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -12.5
$ws.Range("I17").Value = 258
$ws.Range("J17").Value = 244
$ws.Range("K17").Value = 5.737704918032
$ws.Range("L17").Value = 46.590909090909
$ws.Range("M17").Value = 141.121495327103
$ws.Range("N17").Value = -28.729281767955

# Row 18
$ws.Range("D18").NumberFormat = "#,##0"
$ws.Range("D18").Value = 2
$ws.Range("E18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 48
$ws.Range("J18").Value = 50
$ws.Range("K18").Value = -4
$ws.Range("L18").Value = -5.882352941176
$ws.Range("M18").Value = -50
$ws.Range("N18").Value = -91.413237924865

# Row 19
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 15
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = -11.764705882352
$ws.Range("I19").Value = 114
$ws.Range("J19").Value = 118
$ws.Range("K19").Value = -3.389830508474
$ws.Range("L19").Value = -10.9375
$ws.Range("M19").Value = 72.727272727272
$ws.Range("N19").Value = -49.557522123893

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 50
$ws.Range("K20").Value = -48
$ws.Range("L20").Value = -44.680851063829
$ws.Range("M20").Value = -58.064516129032
$ws.Range("N20").Value = -91.900311526479

# Row 21
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -15.789473684210
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = -17.045454545454
$ws.Range("I21").Value = 545
$ws.Range("J21").Value = 573
$ws.Range("K21").Value = -4.886561954624
$ws.Range("L21").Value = 12.603305785124
$ws.Range("M21").Value = 27.336448598130
$ws.Range("N21").Value = -72.165474974463

# Row 22
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = -53.846153846153

# Row 23
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 42.857142857142
$ws.Range("I23").Value = 58
$ws.Range("J23").Value = 60
$ws.Range("K23").Value = -3.333333333333
$ws.Range("L23").Value = -7.936507936507
$ws.Range("M23").Value = 123.076923076923

# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = 27.272727272727
$ws.Range("F24").Value = 73
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 25.862068965517
$ws.Range("I24").Value = 566
$ws.Range("J24").Value = 437
$ws.Range("K24").Value = 29.519450800915
$ws.Range("L24").Value = 42.929292929292
$ws.Range("M24").Value = 113.584905660377

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 1
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 400
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 34
$ws.Range("K25").Value = 185.294117647059
$ws.Range("L25").Value = 36.619718309859

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 14.705882352941
$ws.Range("I26").Value = 339
$ws.Range("J26").Value = 294
$ws.Range("K26").Value = 15.306122448979
$ws.Range("L26").Value = 13
$ws.Range("M26").Value = 19.366197183098

# Row 27
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 6
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 1
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 33
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 37.5
$ws.Range("L27").Value = 50

# Row 28
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 29
$ws.Range("J28").Value = 29
$ws.Range("L28").Value = 0

# Row 29
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("N29").Value = -64.705882352941

# Row 30
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("N30").Value = -68.75

# Row 31
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 2
$ws.Range("K31").Value = -60
$ws.Range("L31").Value = 0
